$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F, matching the style used by the other header cells (B1:E1):
# bold font, thin box border, centered horizontally / top vertically.
$ws.Range("F1").Value = "time_taken"
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$timeTaken = @(
    "2021-10-05 13:40:01.389506",
    "2021-10-05 13:40:01.389519",
    "2021-10-05 13:40:01.389523",
    "2021-10-05 13:40:01.389526",
    "2021-10-05 13:40:01.389530",
    "2021-10-05 13:40:01.389533",
    "2021-10-05 13:40:01.389536",
    "2021-10-05 13:40:01.389539",
    "2021-10-05 13:40:01.389542",
    "2021-10-05 13:40:01.389546",
    "2021-10-05 13:40:01.389549",
    "2021-10-05 13:40:01.389552",
    "2021-10-05 13:40:01.389554",
    "2021-10-05 13:40:01.389557",
    "2021-10-05 13:40:01.389560",
    "2021-10-05 13:40:01.389563",
    "2021-10-05 13:40:01.389567",
    "2021-10-05 13:40:01.389570",
    "2021-10-05 13:40:01.389573",
    "2021-10-05 13:40:01.389576",
    "2021-10-05 13:40:01.389579",
    "2021-10-05 13:40:01.389582",
    "2021-10-05 13:40:01.389585",
    "2021-10-05 13:40:01.389588",
    "2021-10-05 13:40:01.389591",
    "2021-10-05 13:40:01.389594",
    "2021-10-05 13:40:01.389598",
    "2021-10-05 13:40:01.389601",
    "2021-10-05 13:40:01.389604",
    "2021-10-05 13:40:01.389607",
    "2021-10-05 13:40:01.389610",
    "2021-10-05 13:40:01.389612",
    "2021-10-05 13:40:01.389616",
    "2021-10-05 13:40:01.389619",
    "2021-10-05 13:40:01.389622",
    "2021-10-05 13:40:01.389625",
    "2021-10-05 13:40:01.389628",
    "2021-10-05 13:40:01.389630",
    "2021-10-05 13:40:01.389633",
    "2021-10-05 13:40:01.389636",
    "2021-10-05 13:40:01.389640",
    "2021-10-05 13:40:01.389643",
    "2021-10-05 13:40:01.389646",
    "2021-10-05 13:40:01.389649",
    "2021-10-05 13:40:01.389652",
    "2021-10-05 13:40:01.389655",
    "2021-10-05 13:40:01.389658",
    "2021-10-05 13:40:01.389661",
    "2021-10-05 13:40:01.389664",
    "2021-10-05 13:40:01.389667",
    "2021-10-05 13:40:01.389670",
    "2021-10-05 13:40:01.389673",
    "2021-10-05 13:40:01.389676",
    "2021-10-05 13:40:01.389679",
    "2021-10-05 13:40:01.389682",
    "2021-10-05 13:40:01.389685",
    "2021-10-05 13:40:01.389688",
    "2021-10-05 13:40:01.389691",
    "2021-10-05 13:40:01.389694",
    "2021-10-05 13:40:01.389697",
    "2021-10-05 13:40:01.389700",
    "2021-10-05 13:40:01.389703",
    "2021-10-05 13:40:01.389706",
    "2021-10-05 13:40:01.389709",
    "2021-10-05 13:40:01.389714",
    "2021-10-05 13:40:01.389717"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
